$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Roadmap")

# Row 7 (Learn OS basics...): Status In Progress -> Done; In Progress? checkbox cleared; Done? checkbox set
$ws.Range("C7").Value = "Done"
$ws.Range("D7").Value = "☐"
$ws.Range("E7").Value = "☑"

# Row 8 (Learn networking fundamentals...): Status Not Started -> In Progress; In Progress? checkbox set
$ws.Range("C8").Value = "In Progress"
$ws.Range("D8").Value = "☑"

# Update the active selection to C8
$ws.Range("C8").Select() | Out-Null
